$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels in row 1 (A1 and D1 are newly populated;
# B1 and E1 get updated/renamed text)
$ws.Range("A1").Value = "est.RCT"
$ws.Range("B1").Value = "ci.RCT"
$ws.Range("D1").Value = "est.obs"
$ws.Range("E1").Value = "ci.obs"

# Restore selection to E1 (matches final sheetView selection in diff)
$ws.Range("E1").Select()
